{"js": "// Replace the 25 two-digit multiplication problems/answers in the table\n// with the updated set, matching each old expression exactly and\n// substituting the corresponding new expression.\nconst pairs = [\n  [\"56\u00d727=1512\", \"34\u00d739=1326\"],\n  [\"55\u00d759=3245\", \"47\u00d760=2820\"],\n  [\"60\u00d718=1080\", \"96\u00d799=9504\"],\n  [\"63\u00d744=2772\", \"82\u00d791=7462\"],\n  [\"37\u00d761=2257\", \"97\u00d787=8439\"],\n  [\"80\u00d713=1040\", \"51\u00d767=3417\"],\n  [\"36\u00d799=3564\", \"63\u00d762=3906\"],\n  [\"13\u00d720=260\", \"46\u00d744=2024\"],\n  [\"42\u00d760=2520\", \"61\u00d776=4636\"],\n  [\"23\u00d727=621\", \"23\u00d744=1012\"],\n  [\"23\u00d722=506\", \"49\u00d767=3283\"],\n  [\"22\u00d759=1298\", \"71\u00d726=1846\"],\n  [\"55\u00d733=1815\", \"51\u00d740=2040\"],\n  [\"74\u00d753=3922\", \"59\u00d790=5310\"],\n  [\"30\u00d799=2970\", \"38\u00d751=1938\"],\n  [\"77\u00d758=4466\", \"72\u00d772=5184\"],\n  [\"83\u00d716=1328\", \"62\u00d736=2232\"],\n  [\"17\u00d790=1530\", \"21\u00d712=252\"],\n  [\"27\u00d742=1134\", \"17\u00d733=561\"],\n  [\"25\u00d753=1325\", \"80\u00d712=960\"],\n  [\"45\u00d755=2475\", \"54\u00d796=5184\"],\n  [\"76\u00d769=5244\", \"80\u00d745=3600\"],\n  [\"66\u00d741=2706\", \"90\u00d722=1980\"],\n  [\"43\u00d768=2924\", \"38\u00d762=2356\"],\n  [\"44\u00d788=3872\", \"23\u00d724=552\"]\n];\n\nconst body = context.document.body;\nlet totalReplacements = 0;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n    totalReplacements++;\n  }\n  await context.sync();\n}\n\nreturn \"Replaced \" + totalReplacements + \" occurrences across \" + pairs.length + \" pairs.\";\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"56\u00d727=1512\", \"34\u00d739=1326\"),\n  @(\"55\u00d759=3245\", \"47\u00d760=2820\"),\n  @(\"60\u00d718=1080\", \"96\u00d799=9504\"),\n  @(\"63\u00d744=2772\", \"82\u00d791=7462\"),\n  @(\"37\u00d761=2257\", \"97\u00d787=8439\"),\n  @(\"80\u00d713=1040\", \"51\u00d767=3417\"),\n  @(\"36\u00d799=3564\", \"63\u00d762=3906\"),\n  @(\"13\u00d720=260\", \"46\u00d744=2024\"),\n  @(\"42\u00d760=2520\", \"61\u00d776=4636\"),\n  @(\"23\u00d727=621\", \"23\u00d744=1012\"),\n  @(\"23\u00d722=506\", \"49\u00d767=3283\"),\n  @(\"22\u00d759=1298\", \"71\u00d726=1846\"),\n  @(\"55\u00d733=1815\", \"51\u00d740=2040\"),\n  @(\"74\u00d753=3922\", \"59\u00d790=5310\"),\n  @(\"30\u00d799=2970\", \"38\u00d751=1938\"),\n  @(\"77\u00d758=4466\", \"72\u00d772=5184\"),\n  @(\"83\u00d716=1328\", \"62\u00d736=2232\"),\n  @(\"17\u00d790=1530\", \"21\u00d712=252\"),\n  @(\"27\u00d742=1134\", \"17\u00d733=561\"),\n  @(\"25\u00d753=1325\", \"80\u00d712=960\"),\n  @(\"45\u00d755=2475\", \"54\u00d796=5184\"),\n  @(\"76\u00d769=5244\", \"80\u00d745=3600\"),\n  @(\"66\u00d741=2706\", \"90\u00d722=1980\"),\n  @(\"43\u00d768=2924\", \"38\u00d762=2356\"),\n  @(\"44\u00d788=3872\", \"23\u00d724=552\")\n)\n\n$replacedCount = 0\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $found = $find.Execute()\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n  $range.Text = $newText\n  $replacedCount = $replacedCount + 1\n}\n\nWrite-Output \"Replaced $replacedCount occurrences across $($pairs.Length) pairs.\"\n"}
